$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("K2").Value = 0.7887805686561092
$ws.Range("L2").Value = 0.8680469733657697
$ws.Range("M2").Value = 0.886599188681837
$ws.Range("N2").Value = 0.8136557880180363
$ws.Range("O2").Value = 2.421804666519165
$ws.Range("P2").Value = 9.87317681312561

$ws.Range("K3").Value = 0.6406127450008707
$ws.Range("L3").Value = 0.8085239579549567
$ws.Range("M3").Value = 0.8353275031091926
$ws.Range("N3").Value = 0.6820897889154897
$ws.Range("O3").Value = 0.009971857070922852
$ws.Range("P3").Value = 17.88996005058289

$ws.Range("K4").Value = 0.6416683145038612
$ws.Range("L4").Value = 0.8002198692333028
$ws.Range("M4").Value = 0.8282415036296485
$ws.Range("N4").Value = 0.6831840136352011
$ws.Range("O4").Value = 5.420705080032349
$ws.Range("P4").Value = 13.32115316390991

$ws.Range("K5").Value = 0.6474648589640518
$ws.Range("L5").Value = 0.8116937274579437
$ws.Range("M5").Value = 0.8381024182563098
$ws.Range("N5").Value = 0.6884617519931663
$ws.Range("O5").Value = 0.01810789108276367
$ws.Range("P5").Value = 10.53526401519775

$ws.Range("K6").Value = 0.7543611620564029
$ws.Range("L6").Value = 0.8247614623016217
$ws.Range("M6").Value = 0.8492930056684334
$ws.Range("N6").Value = 0.7832130971824748
$ws.Range("O6").Value = 2.326012134552002
$ws.Range("P6").Value = 10.00130295753479

$ws.Range("K7").Value = 0.7866866638067317
$ws.Range("L7").Value = 0.8595642958532245
$ws.Range("M7").Value = 0.8792428206668126
$ws.Range("N7").Value = 0.8117675522933305
$ws.Range("O7").Value = 0.01061010360717773
$ws.Range("P7").Value = 9.8405921459198

$ws.Range("K8").Value = 0.7618932795693877
$ws.Range("L8").Value = 0.8381243637668739
$ws.Range("M8").Value = 0.8607266686378819
$ws.Range("N8").Value = 0.789868640210356
$ws.Range("O8").Value = 2.510090112686157
$ws.Range("P8").Value = 9.885471820831299

$ws.Range("K9").Value = 0.6501465437034426
$ws.Range("L9").Value = 0.8169721950153157
$ws.Range("M9").Value = 0.8425622134110246
$ws.Range("N9").Value = 0.6904993654734219
$ws.Range("O9").Value = 0.01228117942810059
$ws.Range("P9").Value = 19.89185190200806
